# Updates the crypto price/volume columns (D, E) to the latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.812.61"
$ws.Range("E2").Value = "  +1.06%  "
$ws.Range("D3").Value = "3.015.59"
$ws.Range("E3").Value = "  +3.59%  "
$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "'382.46"
$ws.Range("E5").Value = "  +5.37%  "
$ws.Range("D6").Value = "'106.74"
$ws.Range("E6").Value = "  +3.23%  "
$ws.Range("D7").Value = "'0.550"
$ws.Range("E7").Value = "  +1.46%  "
$ws.Range("D8").Value = "'0.998"
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  +2.70%  "
$ws.Range("D10").Value = "'38.06"
$ws.Range("E10").Value = "  +3.49%  "
$ws.Range("E11").Value = "  +0.35%  "
$ws.Range("D12").Value = "'0.0850"
$ws.Range("E12").Value = "  +1.97%  "
$ws.Range("D13").Value = "'18.96"
$ws.Range("E13").Value = "  +2.21%  "
$ws.Range("D14").Value = "3.482.46"
$ws.Range("E14").Value = "  +3.29%  "
$ws.Range("D15").Value = "'7.58"
$ws.Range("E15").Value = "  +3.28%  "
$ws.Range("D16").Value = "2.994.72"
$ws.Range("E16").Value = "  +2.72%  "
$ws.Range("D17").Value = "'0.981"
$ws.Range("E17").Value = "  +2.88%  "
$ws.Range("D18").Value = "51.795.52"
$ws.Range("E18").Value = "  +1.20%  "
$ws.Range("D19").Value = "'3.44"
$ws.Range("E19").Value = "  +4.40%  "
$ws.Range("D20").Value = "'7.50"
$ws.Range("E20").Value = "  +3.97%  "
$ws.Range("D21").Value = "'13.18"
$ws.Range("E21").Value = "  +1.24%  "
$ws.Range("D22").Value = "0.0₃0968"
$ws.Range("E22").Value = "  +2.48%  "
$ws.Range("D23").Value = "'69.11"
$ws.Range("E23").Value = "  +1.17%  "
$ws.Range("D24").Value = "'265.03"
$ws.Range("E24").Value = "  +2.08%  "
$ws.Range("E25").Value = "  +4.22%  "
$ws.Range("E26").Value = "  -0.48%  "
$ws.Range("D27").Value = "'7.31"
$ws.Range("E27").Value = "  +18.47%  "
$ws.Range("D28").Value = "'7.59"
$ws.Range("E28").Value = "  +4.79%  "
$ws.Range("D29").Value = "'26.33"
$ws.Range("E29").Value = "  +1.20%  "
$ws.Range("E30").Value = "  +0.05%  "
$ws.Range("E31").Value = "  -1.57%  "
$ws.Range("D32").Value = "'10.02"
$ws.Range("E32").Value = "  +0.61%  "
$ws.Range("D33").Value = "'35.16"
$ws.Range("E33").Value = "  +0.70%  "
$ws.Range("D34").Value = "'51.20"
$ws.Range("E34").Value = "  +1.15%  "
$ws.Range("E35").Value = "  -2.86%  "
$ws.Range("D36").Value = "'0.0451"
$ws.Range("E36").Value = "  +6.50%  "
$ws.Range("E37").Value = "  -0.04%  "
$ws.Range("D38").Value = "'3.14"
$ws.Range("E38").Value = "  -0.42%  "
$ws.Range("D39").Value = "'17.70"
$ws.Range("E39").Value = "  +4.36%  "
$ws.Range("D40").Value = "'2.67"
$ws.Range("E40").Value = "  -4.44%  "
$ws.Range("D41").Value = "'1.89"
$ws.Range("E41").Value = "  +1.43%  "
$ws.Range("E42").Value = "  +3.59%  "
$ws.Range("D43").Value = "'124.85"
$ws.Range("E43").Value = "  +4.25%  "
$ws.Range("D44").Value = "'22.59"
$ws.Range("E44").Value = "  +0.83%  "
$ws.Range("D45").Value = "'2.09"
$ws.Range("E45").Value = "  -2.76%  "
$ws.Range("E46").Value = "  +7.65%  "
$ws.Range("D47").Value = "'0.277"
$ws.Range("E47").Value = "  +17.02%  "
$ws.Range("D48").Value = "2.064.19"
$ws.Range("E48").Value = "  -0.39%  "
$ws.Range("D49").Value = "'3.31"
$ws.Range("E49").Value = "  +3.76%  "
$ws.Range("D50").Value = "'0.0357"
$ws.Range("E50").Value = "  +14.94%  "
$ws.Range("D51").Value = "'1.32"
$ws.Range("E51").Value = "  +3.58%  "
